# [Resource Manager] Add Performance Field for Host
#
# This script reproduces (as closely as the COM surface allows) the change that:
#  - inserts a new "性能" (Performance) column between 用途(M) and 磁盘(N->O)
#  - sets Compute/Storage/General performance labels per host plus a
#    Medium/High/Low rating in the new column
#  - enriches the disk JSON blobs with a "type" field for the first two
#    hosts, and adds a third host (主机3) row with its own disk JSON
#  - tweaks a handful of column widths / row heights to fit the extra column
#
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new column N ("性能") before the existing 磁盘 column, shifting
#    the disk column (and everything after it) one slot to the right (->O).
# ---------------------------------------------------------------------------
$ws.Columns.Item(14).Insert(-4161)

# ---------------------------------------------------------------------------
# 2. Header row
# ---------------------------------------------------------------------------
$ws.Range("N1").Value2 = "性能"

# ---------------------------------------------------------------------------
# 3. Row 2 (主机1) updates: 用途 -> "Compute", new 性能 -> "Medium",
#    磁盘(now O2) JSON gains "type" fields.
# ---------------------------------------------------------------------------
$ws.Range("M2").Value2 = "Compute"
$ws.Range("N2").Value2 = "Medium"
$ws.Range("O2").Value2 = '[{"name": "vda","capacity": 256,"status": 1, "path": "/", "type": "sata"}, {"name": "vdd", "capacity": 512,"status": 0, "path": "/mnt/vdd", "type":"ssd"}, {"name": "vde","capacity": 1024,"status": 0,"path": "/mnt/vde", "type":"nvme_ssd"}]'

# ---------------------------------------------------------------------------
# 4. Row 3 (主机2) updates: 用途 -> "Storage", new 性能 -> "High",
#    磁盘(now O3) JSON gains "type" fields.
# ---------------------------------------------------------------------------
$ws.Range("M3").Value2 = "Storage"
$ws.Range("N3").Value2 = "High"
$ws.Range("O3").Value2 = '[{"name": "nvme0p1","capacity": 256,"status": 1, "path": "/","type":"nvme_ssd"}, {"name": "nvme0p2", "capacity": 1024,"status": 0, "path": "/mnt/path1", "type": "nvme_ssd"}, {"name": "nvme0p3","capacity": 4096,"status": 0,"path": "/mnt/path2", "type": "nvme_ssd"}]'

# ---------------------------------------------------------------------------
# 5. New row 4 (主机3)
# ---------------------------------------------------------------------------
$ws.Range("A4").Value2 = "主机3"
$ws.Range("B4").Value2 = "192.168.56.13"
$ws.Range("C4").Value2 = "root"
$ws.Range("D4").Value2 = "admin"
$ws.Range("E4").Value2 = "DataCenter1"
$ws.Range("F4").Value2 = "Zone3"
$ws.Range("G4").Value2 = "Rack3"
$ws.Range("H4").Value2 = "Ubuntu"
$ws.Range("I4").Value2 = "5.2.0"
$ws.Range("J4").Value2 = 16
$ws.Range("K4").Value2 = 64
$ws.Range("L4").Value2 = "10GE"
$ws.Range("M4").Value2 = "General"
$ws.Range("N4").Value2 = "Low"
$ws.Range("O4").Value2 = '[{"name": "sda","capacity": 256,"status": 1, "path": "/"}, {"name": "sdb", "capacity": 1024,"status": 0, "path": "/mnt/path1"}, {"name": "sdc","capacity": 4096,"status": 0,"path": "/mnt/path2"}]'

# ---------------------------------------------------------------------------
# 6. Styling: copy the existing "fontId2, vertical-center, no wrap" cell
#    format (already used by C2/D2 in the original workbook) onto every
#    cell that should carry it, then derive the wrap-text variant from it
#    for the disk/JSON column.
# ---------------------------------------------------------------------------
$fmtSrc = $ws.Range("C2")
$fmtSrc.Copy()
$noWrapTargets = $ws.Range("N1,C1:D1,M2:N2,M3:N3,A4:D4,F4:G4,M4:N4")
$noWrapTargets.PasteSpecial(-4122)
$excel.CutCopyMode = 0

$wrapTargets = $ws.Range("O2:O4")
$fmtSrc.Copy()
$wrapTargets.PasteSpecial(-4122)
$wrapTargets.WrapText = $true
$excel.CutCopyMode = 0

Write-Host "stage1 done"
